# Add a new error-message row (code 211, "General" category, "warn" level)
# to the "Error Codes - STIX Elevator" sheet, just below the existing 210
# row and above the 301 row, per commit "added 211 to messages".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 13 (pushes the former row 13..121 down to 14..122)
$ws.Rows.Item(13).Insert()

# Row 13 sits in the same "2xx / General" formatting block as row 12 just
# above it, so copy that row's formatting onto the freshly inserted row.
$ws.Range("A12:F12").Copy()
$ws.Range("A13:F13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new row's content.
$ws.Range("A13").Value = "silent option is not compatible with a policy"
$ws.Range("B13").Value = "General"
$ws.Range("C13").Value = 211
$ws.Range("D13").Value = "warn"
$ws.Range("E13").Value = "ElevatorOptions.__init__()"
# F13 (Notes) stays empty, matching the row's formatting block.

# Match the workbook's resulting selection/cursor position.
[void]$ws.Range("F13").Select()
